$d = $word.ActiveDocument

$replacements = @(
    @{old="2023-12-30 Saturday"; new="2023-12-31 Sunday"},
    @{old="70×28="; new="50×63="},
    @{old="66×14="; new="92×92="},
    @{old="59×60="; new="75×44="},
    @{old="16×14="; new="54×73="},
    @{old="35×88="; new="99×33="},
    @{old="65×54="; new="29×81="},
    @{old="61×23="; new="21×34="},
    @{old="22×50="; new="84×79="},
    @{old="84×61="; new="39×65="},
    @{old="44×60="; new="67×69="},
    @{old="60×73="; new="62×27="},
    @{old="26×34="; new="99×57="},
    @{old="96×12="; new="24×20="},
    @{old="54×55="; new="57×28="},
    @{old="54×45="; new="63×26="},
    @{old="82×81="; new="39×28="},
    @{old="24×52="; new="14×20="},
    @{old="20×25="; new="47×60="},
    @{old="41×89="; new="54×26="},
    @{old="18×41="; new="26×20="},
    @{old="38×99="; new="47×52="},
    @{old="37×94="; new="49×75="},
    @{old="13×69="; new="99×11="},
    @{old="68×15="; new="66×26="},
    @{old="45×29="; new="96×39="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
